$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original style/formatting of column D data range, then force Text
# number format while assigning values so numeric-looking strings (e.g. "1.002")
# are not auto-converted into floating point numbers by Excel.
$dRange = $ws.Range("D2:D51")
$origStyle = $ws.Range("D2").Style
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.961.62"
$ws.Range("E2").Value = "  -5.45%  "
$ws.Range("D3").Value = "1.820.89"
$ws.Range("E3").Value = "  -4.50%  "
$ws.Range("E4").Value = "  -0.94%  "
$ws.Range("D5").Value = "330.04"
$ws.Range("E5").Value = "  -2.43%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").Value = "0.4634"
$ws.Range("E7").Value = "  -2.72%  "
$ws.Range("D8").Value = "0.3844"
$ws.Range("E8").Value = "  -4.11%  "
$ws.Range("D9").Value = "45.86"
$ws.Range("E9").Value = "  -3.86%  "
$ws.Range("D10").Value = "0.07830"
$ws.Range("E10").Value = "  -2.74%  "
$ws.Range("D11").Value = "0.9582"
$ws.Range("E11").Value = "  -3.14%  "
$ws.Range("D12").Value = "21.90"
$ws.Range("E12").Value = "  -5.87%  "
$ws.Range("D13").Value = "1.863.96"
$ws.Range("E13").Value = "  -4.47%  "
$ws.Range("D14").Value = "5.640"
$ws.Range("E14").Value = "  -4.88%  "
$ws.Range("D15").Value = "6.847"
$ws.Range("E15").Value = "  -3.92%  "
$ws.Range("D16").Value = "0.06860"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "86.44"
$ws.Range("E18").Value = "  -3.32%  "
$ws.Range("D19").Value = "0.000009944"
$ws.Range("E19").Value = "  -2.52%  "
$ws.Range("D20").Value = "16.70"
$ws.Range("E20").Value = "  -3.88%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("D22").Value = "27.997.80"
$ws.Range("E22").Value = "  -5.34%  "
$ws.Range("D23").Value = "5.317"
$ws.Range("E23").Value = "  -3.70%  "
$ws.Range("D24").Value = "10.92"
$ws.Range("E24").Value = "  -6.21%  "
$ws.Range("D25").Value = "2.097"
$ws.Range("E25").Value = "  -2.50%  "
$ws.Range("D26").Value = "2.037.66"
$ws.Range("E26").Value = "  -6.61%  "
$ws.Range("D27").Value = "151.67"
$ws.Range("E27").Value = "  -3.33%  "
$ws.Range("D28").Value = "19.18"
$ws.Range("E28").Value = "  -2.61%  "
$ws.Range("D29").Value = "5.707"
$ws.Range("E29").Value = "  -12.05%  "
$ws.Range("D30").Value = "1.963"
$ws.Range("E30").Value = "  -4.42%  "
$ws.Range("D31").Value = "116.33"
$ws.Range("E31").Value = "  -2.60%  "
$ws.Range("D32").Value = "0.9398"
$ws.Range("E32").Value = "  -5.50%  "
$ws.Range("D33").Value = "0.09260"
$ws.Range("E33").Value = "  -2.75%  "
$ws.Range("D34").Value = "5.274"
$ws.Range("E34").Value = "  -4.01%  "
$ws.Range("D35").Value = "3.420"
$ws.Range("E35").Value = "  -3.52%  "
$ws.Range("D36").Value = "1.310"
$ws.Range("E36").Value = "  -5.59%  "
$ws.Range("D37").Value = "0.05944"
$ws.Range("E37").Value = "  -7.96%  "
$ws.Range("D38").Value = "0.02145"
$ws.Range("E38").Value = "  -4.29%  "
$ws.Range("D39").Value = "1.149"
$ws.Range("E39").Value = "  -3.97%  "
$ws.Range("D40").Value = "1.001"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").Value = "7.558"
$ws.Range("E41").Value = "  -2.80%  "
$ws.Range("D42").Value = "0.5581"
$ws.Range("E42").Value = "  -4.15%  "
$ws.Range("D43").Value = "9.885"
$ws.Range("E43").Value = "  -6.40%  "
$ws.Range("D44").Value = "0.1766"
$ws.Range("E44").Value = "  -2.97%  "
$ws.Range("D45").Value = "1.240"
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").Value = "2.225"
$ws.Range("E46").Value = "  -9.38%  "
$ws.Range("D47").Value = "11.59"
$ws.Range("E47").Value = "  -4.18%  "
$ws.Range("E48").Value = "  -4.03%  "
$ws.Range("D49").Value = "0.07002"
$ws.Range("E49").Value = "  -5.70%  "
$ws.Range("D50").Value = "1.824"
$ws.Range("E50").Value = "  -5.90%  "
$ws.Range("D51").Value = "112.42"
$ws.Range("E51").Value = "  -3.28%  "

# Restore the original (default) style now that values are safely stored as text
$dRange.Style = $origStyle
